# "upgrade left table until javakheti"
#
# For the Dmanisi stillbirths sheet:
#   - mark every year value in the "Urban" row (row 6) as confidential/unavailable
#   - mark the 2014 and 2016 values in the "Rural" row (row 7) as confidential/unavailable
#   - rename the worksheet tab from "1" to "Dmanisi"
#   - normalize the confidential-data placeholder from the ellipsis character "…" to
#     three literal dots "..."
#   - remove the stray blank row 8 so the trailing "Note:" row moves up from row 9 to row 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ellipsis = [char]0x2026

# Row 6 ("Urban"): every year column becomes the confidential/unavailable placeholder.
"B6","C6","D6","E6","F6","G6","H6","I6","J6","K6","L6","M6","N6","O6" | ForEach-Object {
    $ws.Range($_).Value = $ellipsis
}

# Row 7 ("Rural"): only the 2014 (F7) and 2016 (H7) values become the placeholder.
"F7","H7" | ForEach-Object {
    $ws.Range($_).Value = $ellipsis
}

# The placeholder text itself changes from the single ellipsis glyph "…" to "...".
# Replacing across the whole sheet updates the shared string in place (and now also
# covers every cell we just set above).
[void]$ws.Cells.Replace($ellipsis, "...")

# Rename the worksheet tab.
$ws.Name = "Dmanisi"

# Row 8 was already empty; deleting it shifts the "Note:" row up from 9 to 8,
# shrinking the used range to A1:X8.
$ws.Rows.Item(8).Delete()
